$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 (was: ECs / Dkk2 / Kremen2 / ECs) -> becomes: FAPs / Dkk2 / Kremen2 / MuSCs ---
$ws.Range("A2").Value = "FAPs"
$ws.Range("B2").Value = "Dkk2"
$ws.Range("C2").Value = "Kremen2"
$ws.Range("D2").Value = "MuSCs"

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 4.066751
$ws.Range("H2").Value = 12.200253
$ws.Range("I2").Value = 0.9827679245700983
$ws.Range("J2").Value = 0.9827679245700983
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.5
$ws.Range("M2").Value = 0.08043649999999999
$ws.Range("N2").Value = 0.160873
$ws.Range("O2").Value = 1
$ws.Range("P2").Value = 1
$ws.Range("Q2").Value = 0.3271152168115
$ws.Range("R2").Value = 1.962691300869
$ws.Range("S2").Value = 0.9827679245700983
$ws.Range("T2").Value = 0.9827679245700983

# --- Row 3 (was: FAPs / Dkk2 / Kremen2 / ECs) -> becomes: Resolving-Mac / Dkk2 / Kremen2 / MuSCs ---
$ws.Range("A3").Value = "Resolving-Mac"
$ws.Range("B3").Value = "Dkk2"
$ws.Range("C3").Value = "Kremen2"
$ws.Range("D3").Value = "MuSCs"

$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.07130733333333333
$ws.Range("H3").Value = 0.213922
$ws.Range("I3").Value = 0.0172320754299017
$ws.Range("J3").Value = 0.0172320754299017
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 0.5
$ws.Range("M3").Value = 0.08043649999999999
$ws.Range("N3").Value = 0.160873
$ws.Range("O3").Value = 1
$ws.Range("P3").Value = 1
$ws.Range("Q3").Value = 0.005735712317666666
$ws.Range("R3").Value = 0.034414273906
$ws.Range("S3").Value = 0.0172320754299017
$ws.Range("T3").Value = 0.0172320754299017

# --- Row 4 (was: MuSCs / Dkk2 / Kremen2 / ECs) is removed entirely ---
$ws.Range("A4").EntireRow.Delete()
